$d = $word.ActiveDocument

# The "Referências" section has two citations shaped like:
#   ... Disponível em: <URL>. Acesso em 16 Jan. 2021.
# where the angle brackets around the hyperlinked URL are their own
# standalone runs. The edit drops those bracket runs so the citation
# reads:
#   ... Disponível em: URL. Acesso em 16 Jan. 2021.
#
# Each ">" (or "> ") run sits immediately after a w:hyperlink field and
# shares identical run formatting with the plain-text run that follows
# it (". Acesso em 16 "); deleting the "<" run first would let the host
# engine's run-renormalization merge those two together (changing the
# document in a way the diff doesn't call for). So remove the ">" /"> "
# runs before touching the "<" runs, keeping every other run intact.

# 1) First citation (QEDU): "...taxas-rendimento>. Acesso..." -> delete the lone ">" run.
$r = $d.Content
if (-not $r.Find.Execute(">")) { throw "could not find first '>' run" }
$r.Delete()

# 2) Second citation (INEP): "...escolar.pdf> . Acesso..." -> delete the "> " run (char + trailing space).
$r = $d.Content
if (-not $r.Find.Execute("> ")) { throw "could not find '> ' run" }
$r.Delete()

# 3) First citation: "Disponível em: <https://www.qedu..." -> delete the lone "<" run.
$r = $d.Content
if (-not $r.Find.Execute("<")) { throw "could not find first '<' run" }
$r.Delete()

# 4) Second citation: "Disponível em: <https://download.inep..." -> delete the lone "<" run.
$r = $d.Content
if (-not $r.Find.Execute("<")) { throw "could not find second '<' run" }
$r.Delete()
